$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): split "document"/"git" labels into
# "document status" / "git status", add a new "work done by" column,
# and rename the first header from "Java Tutorial" to "Concept".
$ws.Range("B1").Value = "work done by"
$ws.Range("C1").Value = "document status"
$ws.Range("D1").Value = "git status"
$ws.Range("A1").Value = "Concept"

# Widen column B (new "work done by" column) and column C
# ("document status" column).
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(3).ColumnWidth = 17
